# "Generate Report for handoff"
#
# The localization status workbook tracks per-file handoff state across three
# sheets: "Overview" (sheet1), "zh-cn" (sheet2) and "de-de" (sheet3).
#
# This run represents a new handoff cycle:
#   - the source file that used to be f711b5db-...md was renamed/replaced by
#     bb85d3f6-...md (still "Ready for handoff", still produced a fresh .xlf)
#   - a brand-new source file 6286a789-...md showed up whose handoff
#     transform failed ("Handoff transform failed" / "Ignored")
#   - the always-present ".localization-config" bookkeeping row shifts down
#     to make room for the new row
#
# We rebuild each sheet's rows from scratch (after clearing old hyperlinks)
# so that shared strings / hyperlink relationships never end up stale.

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/022543b1b52e19e65ba5288fe19c2e3b849e31cb/e2e/"
$cfgUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/022543b1b52e19e65ba5288fe19c2e3b849e31cb/.localization-config"
$zhXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3cc08f43afa966d3353e56c93e6e73573cd31e34/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/"
$deXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bb35ec4cdf46755634d69239aa6dd11d69db34cc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/"

$oldMdName = "f711b5db-7719-4c11-b8a5-f80bec1ef276.md"
$newGuid = "bb85d3f6-00ee-4245-8fdf-c04ac291d568"
$newMdName = "$newGuid.md"
$failedMdName = "6286a789-0224-4352-b49a-cd3226406aca.md"
$cfgName = ".localization-config"

$oldHash = "3cb59742098c517f4e4997295045150eafc24da9"
$newHash = "38a774548b0d1ceca33eed15e97804f4b11c5d46"

$zhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$deXlfName = "$newGuid.$newHash.de-de.xlf"

$zhHandoffDatetime = "2016-01-18 04:03:34"
$deHandoffDatetime = "2016-01-18 04:03:45"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

$ws1.Range("A1").Value = "File Name"
$ws1.Range("B1").Value = "zh-cn"
$ws1.Range("C1").Value = "de-de"

$ws1.Range("A2").Value = $newMdName
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = $failedMdName
$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

$ws1.Range("A4").Value = $cfgName
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), ($mdBase + $newMdName), "", "", $newMdName) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($mdBase + $failedMdName), "", "", $failedMdName) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), $cfgUrl, "", "", $cfgName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Range("A1").Value = "Source File Name"
$ws2.Range("B1").Value = "Status"
$ws2.Range("C1").Value = "Latest Handoff File"
$ws2.Range("D1").Value = "Latest Handoff Datetime"
$ws2.Range("E1").Value = "Latest Target File"
$ws2.Range("F1").Value = "Latest Handback File"
$ws2.Range("G1").Value = "Latest Handback DateTime"
$ws2.Range("H1").Value = "Handoff Reason"
$ws2.Range("I1").Value = "Dependency From"

$ws2.Range("A2").Value = $newMdName
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = $zhXlfName
$ws2.Range("D2").Value = $zhHandoffDatetime
$ws2.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = $failedMdName
$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = $epoch
$ws2.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = "Ignored"

$ws2.Range("A4").Value = $cfgName
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = $epoch
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), ($mdBase + $newMdName), "", "", $newMdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), ($zhXlfBase + $zhXlfName), "", "", $zhXlfName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($mdBase + $failedMdName), "", "", $failedMdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), $cfgUrl, "", "", $cfgName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Range("A1").Value = "Source File Name"
$ws3.Range("B1").Value = "Status"
$ws3.Range("C1").Value = "Latest Handoff File"
$ws3.Range("D1").Value = "Latest Handoff Datetime"
$ws3.Range("E1").Value = "Latest Target File"
$ws3.Range("F1").Value = "Latest Handback File"
$ws3.Range("G1").Value = "Latest Handback DateTime"
$ws3.Range("H1").Value = "Handoff Reason"
$ws3.Range("I1").Value = "Dependency From"

$ws3.Range("A2").Value = $newMdName
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = $deXlfName
$ws3.Range("D2").Value = $deHandoffDatetime
$ws3.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G2").Value = $epoch
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = $failedMdName
$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = $epoch
$ws3.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = "Ignored"

$ws3.Range("A4").Value = $cfgName
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = $epoch
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), ($mdBase + $newMdName), "", "", $newMdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), ($deXlfBase + $deXlfName), "", "", $deXlfName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($mdBase + $failedMdName), "", "", $failedMdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), $cfgUrl, "", "", $cfgName) | Out-Null
